$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.413.59'
$ws.Range('E2').Value = '  -6.19%  '
$ws.Range('D3').Value = '3.128.56'
$ws.Range('E3').Value = '  -7.91%  '
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.56'
$ws.Range('E5').Value = '  -4.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.39'
$ws.Range('E6').Value = '  -10.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  -5.25%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '3.129.04'
$ws.Range('E9').Value = '  -7.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.583'
$ws.Range('E10').Value = '  -7.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '51.09'
$ws.Range('E11').Value = '  -14.07%  '
$ws.Range('E12').Value = '  -6.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000244'
$ws.Range('E13').Value = '  -4.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.62'
$ws.Range('E14').Value = '  -7.61%  '
$ws.Range('D15').Value = '3.638.75'
$ws.Range('E15').Value = '  -8.18%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.112'
$ws.Range('E16').Value = '  -8.97%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.141.35'
$ws.Range('E17').Value = '  -7.96%  '
$ws.Range('D18').Value = '61.371.95'
$ws.Range('E18').Value = '  -6.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.65'
$ws.Range('E19').Value = '  -5.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.56'
$ws.Range('E20').Value = '  -6.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.930'
$ws.Range('E21').Value = '  -5.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '354.68'
$ws.Range('E22').Value = '  -5.41%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '78.45'
$ws.Range('E23').Value = '  -4.58%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.60'
$ws.Range('E24').Value = '  -4.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.65'
$ws.Range('E25').Value = '  -2.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.08'
$ws.Range('E26').Value = '  +3.27%  '
$ws.Range('E27').Value = '  +1.41%  '
$ws.Range('E28').Value = '  -5.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.80'
$ws.Range('E29').Value = '  -7.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.87'
$ws.Range('E30').Value = '  -8.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '631.43'
$ws.Range('E31').Value = '  -9.26%  '
$ws.Range('E32').Value = '  -8.02%  '
$ws.Range('E33').Value = '  -8.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.88'
$ws.Range('E34').Value = '  -3.54%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  -6.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.85'
$ws.Range('E37').Value = '  -8.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.15'
$ws.Range('E38').Value = '  -4.46%  '
$ws.Range('E39').Value = '  -6.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('D41').Value = '0.0₃0671'
$ws.Range('E41').Value = '  +7.08%  '
$ws.Range('E42').Value = '  -6.39%  '
$ws.Range('D43').Value = '2.781.93'
$ws.Range('E43').Value = '  -3.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.40'
$ws.Range('E44').Value = '  +1.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.58'
$ws.Range('E45').Value = '  -3.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.74'
$ws.Range('E46').Value = '  +3.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0373'
$ws.Range('E47').Value = '  -6.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.47'
$ws.Range('E48').Value = '  -11.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.86'
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '131.50'
$ws.Range('E50').Value = '  -4.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.119'
$ws.Range('E51').Value = '  -5.85%  '
